# DB_Films.xlsx - "modif struct-info about DB"
#
# Adds two new lookup sheets describing the users / user-film relation:
#   - "ListFilms" (IdFilms, IdUser)  -> created first
#   - "Users"     (Id, UserName)     -> created second, then moved before ListFilms
# and updates the view state (zoom + selection) of the "Films" sheet.

$wb = $excel.ActiveWorkbook

$wsFilms = $wb.Worksheets.Item("Films")

# --- Create the two new sheets -------------------------------------------------
# "ListFilms" is created first, so it receives the lower internal sheetId.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsListFilms = $wb.Worksheets.Add($null, $lastSheet)
$wsListFilms.Name = "ListFilms"

# "Users" is created second (appended after "ListFilms"), receiving the higher sheetId.
$wsUsers = $wb.Worksheets.Add($null, $wsListFilms)
$wsUsers.Name = "Users"

# Reorder: "Users" tab moves in front of "ListFilms".
$wsUsers.Move($wsListFilms)

# Re-resolve both sheets by name after the Move so later calls bind to the
# correct (post-move) sheet objects.
$wsUsers = $wb.Worksheets.Item("Users")
$wsListFilms = $wb.Worksheets.Item("ListFilms")

# --- Fill in the header rows, reusing the existing "good/green" header style ---
# (copy format only from the Films header so no new cell style is created)
$null = $wsFilms.Range("A1").Copy()
$null = $wsListFilms.Range("A1:B1").PasteSpecial(-4122)
$wsListFilms.Range("A1").Value = "IdFilms"
$wsListFilms.Range("B1").Value = "IdUser"

$null = $wsFilms.Range("A1").Copy()
$null = $wsUsers.Range("A1:B1").PasteSpecial(-4122)
$wsUsers.Range("A1").Value = "Id"
$wsUsers.Range("B1").Value = "UserName"

# --- View / selection bookkeeping ---------------------------------------------
# Select every cell on "ListFilms" ...
$null = $wsListFilms.Cells.Select()

# ... then land on "Users" (last-active tab) with B2 selected.
$null = $wsUsers.Activate()
$null = $wsUsers.Range("B2").Select()

# Zoom out the "Films" sheet and leave it fully selected (no longer the active tab).
$null = $wsFilms.Select()
$excel.ActiveWindow.Zoom = 85
$null = $wsFilms.Cells.Select()

# Return focus to "Users", matching the final saved state of the workbook.
$null = $wsUsers.Activate()
$null = $wsUsers.Range("B2").Select()
